$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.719.15"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.601.45"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.61"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0619"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.59"
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0841"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "1.825.67"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.604.99"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.524"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.42"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("D17").Value = "26.690.10"
$ws.Range("E18").Value = "  +3.13%  "
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "209.56"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.40"
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.12"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.36"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0523"
$ws.Range("E30").Value = "  +3.10%  "
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.26"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("D34").Value = "1.295.53"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.626"
$ws.Range("E35").Value = "  -4.76%  "
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.09"
$ws.Range("E39").Value = "  +20.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.825"
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.44"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.784"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.25"
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("D45").Value = "1.737.31"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.35"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("E51").Value = "  +0.18%  "
